$wb = $excel.ActiveWorkbook

$wsGame = $wb.Worksheets.Item(1)
$wsLevel1 = $wb.Worksheets.Item(2)
$wsLevel2 = $wb.Worksheets.Item(3)

# --- Level 1: swap the "code" / "sector number" columns -------------------
# Previously column A held the code text and column B held the sector
# number. Now column A holds the sector number and column B holds the code,
# so that a sector number can repeat across several rows (multiple codes in
# one sector).
for ($row = 6; $row -le 10; $row++) {
    $codeCell = $wsLevel1.Cells.Item($row, 1)
    $numCell = $wsLevel1.Cells.Item($row, 2)

    $code = $codeCell.Value()
    $num = $numCell.Value()

    $codeCell.Value = $num
    $codeCell.HorizontalAlignment = 1

    $numCell.Value = $code
}

# --- Level 2: same column swap, plus extra codes for sector 1 -------------
for ($row = 6; $row -le 10; $row++) {
    $codeCell = $wsLevel2.Cells.Item($row, 1)
    $numCell = $wsLevel2.Cells.Item($row, 2)

    $code = $codeCell.Value()
    $num = $numCell.Value()

    $codeCell.Value = $num
    $codeCell.HorizontalAlignment = 1

    $numCell.Value = $code
}

# Sector 1 (row 6) now has two more codes recorded in columns C..E.
$wsLevel2.Range("C6").Value = "code1.1"
$wsLevel2.Range("D6").Value = "code1.2"
$wsLevel2.Range("E6").Value = "code1.3"

# --- Move the "codes start at row 6" comment from col A to col B ----------
$comment1 = $wsLevel1.Range("A6").Comment
$commentText1 = $comment1.Text()
$comment1.Delete()
$wsLevel1.Range("B6").AddComment($commentText1)

$comment2 = $wsLevel2.Range("A6").Comment
$commentText2 = $comment2.Text()
$comment2.Delete()
$wsLevel2.Range("B6").AddComment($commentText2)

# --- Selection / active sheet bookkeeping ----------------------------------
# Leave the "Game" sheet's own selection (B4) untouched.
$wsLevel1.Range("A11").Select() | Out-Null
$wsLevel2.Range("E7").Select() | Out-Null
$wsLevel2.Activate() | Out-Null
